# Fix issue with no color on measures.
# Release 1.1.1
#
# A stray "#00ff00" literal had been entered in the _MasterItemColor cell (E3)
# of the Measure1 row. Clear it so the row has no (incorrect) color set,
# matching the other rows that rely on their own color values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").ClearContents()

# Update the active selection left over from editing that cell.
$ws.Range("E6").Select()
